$wb = $excel.ActiveWorkbook

# --- Sheet1 (GW_PC_AccountCreation): selection change only ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2:F2").Select()

# --- Sheet3 (GW_PC_GoogleSearch): selection change only ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("G12").Select()

# --- Add new sheet GW_BC_BillingSummaryAPI after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "GW_BC_BillingSummaryAPI"

# Cell values are written in the same order the source workbook's
# automation authored them (preserves the original shared-string table
# ordering): B1,D1,C1,A1,A2,E1,F1 first, then H1,I1,G1,G2,H2, then the
# billing-summary columns K1..Z1, then J1, then the Jira ids D2,C2,B2,
# then the remaining K2.. values.
$ws4.Range("B1").Value = "Feature"
$ws4.Range("D1").Value = "Scenario_ID"
$ws4.Range("C1").Value = "TestExecution_ID"
$ws4.Range("A1").Value = "Sprint"
$ws4.Range("A2").Value = "Sprint1"
$ws4.Range("E1").Value = "Iteration"
$ws4.Range("F1").Value = "SubIteration"

$ws4.Range("H1").Value = "PolicyNumber"
$ws4.Range("I1").Value = "TermNumber"
$ws4.Range("G1").Value = "EndPoint"
$ws4.Range("G2").Value = "http://10.109.11.67:8580/bc/ws/gw/webservice/policycenter/bc900/BillingSummaryAPI"
$ws4.Range("H2").Value = "'2459765753"

$ws4.Range("K1").Value = "BillingStatus_BillingMethodCode"
$ws4.Range("L1").Value = "BillingStatus_Delinquent"
$ws4.Range("M1").Value = "BillingStatus_PastDue"
$ws4.Range("N1").Value = "BillingStatus_TotalBilled"
$ws4.Range("O1").Value = "BillingStatus_Unbilled"
$ws4.Range("P1").Value = "CurrentOutstanding"
$ws4.Range("Q1").Value = "Invoices"
$ws4.Range("R1").Value = "Paid"
$ws4.Range("S1").Value = "PaymentPlanName"
$ws4.Range("T1").Value = "PolicyTermInfos_EffectiveDate"
$ws4.Range("U1").Value = "PolicyTermInfos_ExpirationDate"
$ws4.Range("V1").Value = "PolicyTermInfos_PolicyNumber"
$ws4.Range("W1").Value = "PolicyTermInfos_TermNumber"
$ws4.Range("X1").Value = "Retrieved"
$ws4.Range("Y1").Value = "TotalCharges"
$ws4.Range("Z1").Value = "WrittenOff"
$ws4.Range("J1").Value = "Archived"

$ws4.Range("D2").Value = "'@DEMO-6"
$ws4.Range("C2").Value = "'@DEMO-7"
$ws4.Range("B2").Value = "'@DEMO-8"

$ws4.Range("E2").Value = 1
$ws4.Range("F2").Value = 1
$ws4.Range("I2").Value = 1
$ws4.Range("J2").Value = $false
$ws4.Range("K2").Value = "DirectBill"
$ws4.Range("L2").Value = $false
$ws4.Range("M2").Value = "0 usd"
$ws4.Range("N2").Value = "0.00 usd"
$ws4.Range("O2").Value = "0.00 usd"
$ws4.Range("P2").Value = "0 usd"
$ws4.Range("Q2").Value = 0
$ws4.Range("R2").Value = "0.00 usd"
$ws4.Range("S2").Value = "QA1PAYMENTPLAN10"
$ws4.Range("T2").Value = "2017-11-20T00:00:00+05:30"
$ws4.Range("U2").Value = "2018-11-20T00:00:00+05:30"
$ws4.Range("V2").Formula = "=H2"
$ws4.Range("W2").Formula = "=I2"
$ws4.Range("X2").Value = $false
$ws4.Range("Y2").Value = "0.00 usd"
$ws4.Range("Z2").Value = "0.00 usd"

# Hyperlink on the EndPoint data cell
$ws4.Hyperlinks.Add($ws4.Range("G2"), "http://10.109.11.67:8580/bc/ws/gw/webservice/policycenter/bc900/BillingSummaryAPI")

# Column widths (best-fit sizes matching the authored workbook)
$ws4.Columns.Item(1).ColumnWidth = 5.7265625
$ws4.Columns.Item(2).ColumnWidth = 9.6328125
$ws4.Columns.Item(3).ColumnWidth = 15.26953125
$ws4.Columns.Item(4).ColumnWidth = 10.7265625
$ws4.Columns.Item(5).ColumnWidth = 8
$ws4.Columns.Item(6).ColumnWidth = 11.08984375
$ws4.Columns.Item(7).ColumnWidth = 75.1796875
$ws4.Columns.Item(8).ColumnWidth = 12.36328125
$ws4.Columns.Item(9).ColumnWidth = 12
$ws4.Columns.Item(10).ColumnWidth = 12
$ws4.Columns.Item(11).ColumnWidth = 28
$ws4.Columns.Item(12).ColumnWidth = 21.26953125
$ws4.Columns.Item(13).ColumnWidth = 18.81640625
$ws4.Columns.Item(14).ColumnWidth = 20.7265625
$ws4.Columns.Item(15).ColumnWidth = 18.90625
$ws4.Columns.Item(16).ColumnWidth = 17.54296875
$ws4.Columns.Item(17).ColumnWidth = 7.54296875
$ws4.Columns.Item(18).ColumnWidth = 7.81640625
$ws4.Columns.Item(19).ColumnWidth = 19
$ws4.Columns.Item(20).ColumnWidth = 26.453125
$ws4.Columns.Item(21).ColumnWidth = 27.81640625
$ws4.Columns.Item(22).ColumnWidth = 27
$ws4.Columns.Item(23).ColumnWidth = 26.6328125
$ws4.Columns.Item(25).ColumnWidth = 11.6328125
$ws4.Columns.Item(26).ColumnWidth = 9.81640625

# Selection / view state for the new sheet
$ws4.Range("Q3").Select()
